# Apply updated values to result_data_KNN.xlsx (Sheet1)
# Commit message: "Update Name of Algo"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.108
$ws.Range("A3").Value = -21.666
$ws.Range("C3").Value = -13.204
$ws.Range("D6").Value = -7.803999999999999
$ws.Range("C12").Value = -11.288
$ws.Range("A14").Value = -21.603
$ws.Range("D19").Value = -7.860000000000001
$ws.Range("A21").Value = -20.257
$ws.Range("A23").Value = -20.955
$ws.Range("C24").Value = -12.568
$ws.Range("D24").Value = -7.562
$ws.Range("A25").Value = -20.372
$ws.Range("B25").Value = 7.273999999999999
$ws.Range("C25").Value = -12.584
$ws.Range("A26").Value = -21.409
$ws.Range("B27").Value = 5.513
$ws.Range("A29").Value = -21.267
$ws.Range("D30").Value = -7.294
$ws.Range("B31").Value = 5.988
$ws.Range("D31").Value = -8.127000000000001
$ws.Range("D33").Value = -7.622
$ws.Range("B39").Value = 7.5
$ws.Range("D42").Value = -8.395999999999999
$ws.Range("B48").Value = 5.156
$ws.Range("C50").Value = -13.133
$ws.Range("B51").Value = 5.42
$ws.Range("B52").Value = 5.358000000000001
$ws.Range("A53").Value = -21.642
$ws.Range("C53").Value = -11.646
$ws.Range("B55").Value = 4.508
$ws.Range("D55").Value = -8.210000000000001
$ws.Range("B56").Value = 5.169
$ws.Range("A57").Value = -21.322
$ws.Range("B57").Value = 7.019
$ws.Range("C57").Value = -12.872
$ws.Range("D58").Value = -7.970999999999999
$ws.Range("A59").Value = -22.097
$ws.Range("C61").Value = -13.025
$ws.Range("C63").Value = -11.591
$ws.Range("D65").Value = -7.869
$ws.Range("A69").Value = -21.546
$ws.Range("C70").Value = -12.188
$ws.Range("D70").Value = -7.571
$ws.Range("B73").Value = 6.439
$ws.Range("D75").Value = -7.645
$ws.Range("A79").Value = -21.215
$ws.Range("A83").Value = -22.185
$ws.Range("D83").Value = -8.397
$ws.Range("C86").Value = -12.416
$ws.Range("D86").Value = -7.65
$ws.Range("B89").Value = 6.009
$ws.Range("B90").Value = 5.515
$ws.Range("A91").Value = -21.533
$ws.Range("B92").Value = 6.027
$ws.Range("A93").Value = -21.439
$ws.Range("D96").Value = -7.433999999999999
$ws.Range("D97").Value = -8.17
$ws.Range("C98").Value = -12.45
$ws.Range("C100").Value = -12.661
$ws.Range("C102").Value = -13.361
